$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Adam2"
$ws.Cells.Item(2, 3).Value = "Itgb7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1251886666666667
$ws.Cells.Item(2, 8).Value = 0.375566
$ws.Cells.Item(2, 9).Value = 0.6155719715657366
$ws.Cells.Item(2, 10).Value = 0.7060466830097307
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.887307
$ws.Cells.Item(2, 14).Value = 2.661921
$ws.Cells.Item(2, 15).Value = 0.03991909470044044
$ws.Cells.Item(2, 16).Value = 0.04024932703229714
$ws.Cells.Item(2, 17).Value = 0.111080780254
$ws.Cells.Item(2, 18).Value = 0.9997270222860001
$ws.Cells.Item(2, 19).Value = 0.02457307582786947
$ws.Cells.Item(2, 20).Value = 0.02841790384452728

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Adam2"
$ws.Cells.Item(3, 3).Value = "Itgb7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1251886666666667
$ws.Cells.Item(3, 8).Value = 0.375566
$ws.Cells.Item(3, 9).Value = 0.6155719715657366
$ws.Cells.Item(3, 10).Value = 0.7060466830097307
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.9845453333333333
$ws.Cells.Item(3, 14).Value = 2.953636
$ws.Cells.Item(3, 15).Value = 0.04429375447078636
$ws.Cells.Item(3, 16).Value = 0.04466017635322986
$ws.Cells.Item(3, 17).Value = 0.1232539175528889
$ws.Cells.Item(3, 18).Value = 1.109285257976
$ws.Cells.Item(3, 19).Value = 0.02726599376763062
$ws.Cells.Item(3, 20).Value = 0.03153216937682756

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Adam2"
$ws.Cells.Item(4, 3).Value = "Itgb7"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1251886666666667
$ws.Cells.Item(4, 8).Value = 0.375566
$ws.Cells.Item(4, 9).Value = 0.6155719715657366
$ws.Cells.Item(4, 10).Value = 0.7060466830097307
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 15.43767133333333
$ws.Cells.Item(4, 14).Value = 46.313014
$ws.Cells.Item(4, 15).Value = 0.6945260928963797
$ws.Cells.Item(4, 16).Value = 0.7002715882016618
$ws.Cells.Item(4, 17).Value = 1.932621490658222
$ws.Cells.Item(4, 18).Value = 17.393593415924
$ws.Cells.Item(4, 19).Value = 0.4275307963080724
$ws.Cells.Item(4, 20).Value = 0.4944244320557394

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Adam2"
$ws.Cells.Item(5, 3).Value = "Itgb7"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1251886666666667
$ws.Cells.Item(5, 8).Value = 0.375566
$ws.Cells.Item(5, 9).Value = 0.6155719715657366
$ws.Cells.Item(5, 10).Value = 0.7060466830097307
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.547111
$ws.Cells.Item(5, 14).Value = 1.094222
$ws.Cells.Item(5, 15).Value = 0.02461400149063703
$ws.Cells.Item(5, 16).Value = 0.01654508121162658
$ws.Cells.Item(5, 17).Value = 0.06849209660866667
$ws.Cells.Item(5, 18).Value = 0.410952579652
$ws.Cells.Item(5, 19).Value = 0.01515168942571342
$ws.Cells.Item(5, 20).Value = 0.01168159970959556

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Adam2"
$ws.Cells.Item(6, 3).Value = "Itgb7"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.1251886666666667
$ws.Cells.Item(6, 8).Value = 0.375566
$ws.Cells.Item(6, 9).Value = 0.6155719715657366
$ws.Cells.Item(6, 10).Value = 0.7060466830097307
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.370998666666666
$ws.Cells.Item(6, 14).Value = 13.112996
$ws.Cells.Item(6, 15).Value = 0.1966470564417564
$ws.Cells.Item(6, 16).Value = 0.1982738272011845
$ws.Cells.Item(6, 17).Value = 0.5471994950817777
$ws.Cells.Item(6, 18).Value = 4.924795455736
$ws.Cells.Item(6, 19).Value = 0.1210504162364507
$ws.Cells.Item(6, 20).Value = 0.1399905780230408

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Adam2"
$ws.Cells.Item(7, 3).Value = "Itgb7"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.5
$ws.Cells.Item(7, 7).Value = 0.078181
$ws.Cells.Item(7, 8).Value = 0.156362
$ws.Cells.Item(7, 9).Value = 0.3844280284342634
$ws.Cells.Item(7, 10).Value = 0.2939533169902694
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.887307
$ws.Cells.Item(7, 14).Value = 2.661921
$ws.Cells.Item(7, 15).Value = 0.03991909470044044
$ws.Cells.Item(7, 16).Value = 0.04024932703229714
$ws.Cells.Item(7, 17).Value = 0.06937054856699999
$ws.Cells.Item(7, 18).Value = 0.416223291402
$ws.Cells.Item(7, 19).Value = 0.01534601887257097
$ws.Cells.Item(7, 20).Value = 0.01183142318776986

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Adam2"
$ws.Cells.Item(8, 3).Value = "Itgb7"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.5
$ws.Cells.Item(8, 7).Value = 0.078181
$ws.Cells.Item(8, 8).Value = 0.156362
$ws.Cells.Item(8, 9).Value = 0.3844280284342634
$ws.Cells.Item(8, 10).Value = 0.2939533169902694
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.9845453333333333
$ws.Cells.Item(8, 14).Value = 2.953636
$ws.Cells.Item(8, 15).Value = 0.04429375447078636
$ws.Cells.Item(8, 16).Value = 0.04466017635322986
$ws.Cells.Item(8, 17).Value = 0.07697273870533333
$ws.Cells.Item(8, 18).Value = 0.461836432232
$ws.Cells.Item(8, 19).Value = 0.01702776070315574
$ws.Cells.Item(8, 20).Value = 0.01312800697640231

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Adam2"
$ws.Cells.Item(9, 3).Value = "Itgb7"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.5
$ws.Cells.Item(9, 7).Value = 0.078181
$ws.Cells.Item(9, 8).Value = 0.156362
$ws.Cells.Item(9, 9).Value = 0.3844280284342634
$ws.Cells.Item(9, 10).Value = 0.2939533169902694
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 15.43767133333333
$ws.Cells.Item(9, 14).Value = 46.313014
$ws.Cells.Item(9, 15).Value = 0.6945260928963797
$ws.Cells.Item(9, 16).Value = 0.7002715882016618
$ws.Cells.Item(9, 17).Value = 1.206932582511333
$ws.Cells.Item(9, 18).Value = 7.241595495067999
$ws.Cells.Item(9, 19).Value = 0.2669952965883073
$ws.Cells.Item(9, 20).Value = 0.2058471561459225

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Adam2"
$ws.Cells.Item(10, 3).Value = "Itgb7"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.5
$ws.Cells.Item(10, 7).Value = 0.078181
$ws.Cells.Item(10, 8).Value = 0.156362
$ws.Cells.Item(10, 9).Value = 0.3844280284342634
$ws.Cells.Item(10, 10).Value = 0.2939533169902694
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.547111
$ws.Cells.Item(10, 14).Value = 1.094222
$ws.Cells.Item(10, 15).Value = 0.02461400149063703
$ws.Cells.Item(10, 16).Value = 0.01654508121162658
$ws.Cells.Item(10, 17).Value = 0.042773685091
$ws.Cells.Item(10, 18).Value = 0.171094740364
$ws.Cells.Item(10, 19).Value = 0.009462312064923614
$ws.Cells.Item(10, 20).Value = 0.004863481502031019

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Adam2"
$ws.Cells.Item(11, 3).Value = "Itgb7"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.5
$ws.Cells.Item(11, 7).Value = 0.078181
$ws.Cells.Item(11, 8).Value = 0.156362
$ws.Cells.Item(11, 9).Value = 0.3844280284342634
$ws.Cells.Item(11, 10).Value = 0.2939533169902694
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 4.370998666666666
$ws.Cells.Item(11, 14).Value = 13.112996
$ws.Cells.Item(11, 15).Value = 0.1966470564417564
$ws.Cells.Item(11, 16).Value = 0.1982738272011845
$ws.Cells.Item(11, 17).Value = 0.3417290467586666
$ws.Cells.Item(11, 18).Value = 2.050374280552
$ws.Cells.Item(11, 19).Value = 0.07559664020530572
$ws.Cells.Item(11, 20).Value = 0.05828324917814368
